$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sort the data range (A2:D18) in ascending order by column A (time).
$rng = $ws.Range("A2:D18")
$key1 = $ws.Range("A2:A18")
$rng.Sort($key1, 1)
